# Weekly "Fruta / hortaliza" update: a new week's worth of data (rows 93-94)
# is inserted at the top of this category's data block, the rest of the
# weekly rows shift down by one week (2 rows), and the two rows that fall
# off the bottom (old rows 215-216) are appended as new rows 217-218.
#
# Only the Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M), Precio $/Kg (P) and Fecha (D) columns carry
# week-specific data; every other column (A,B,C,E,F,G,H,I,N,O,Q,R) only
# depends on whether the row is a "Primera" (odd) or "Segunda" (even) row
# and is left untouched by the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 93
$lastDataRow = 216
$newLastRow = 218

# Columns (1-based) that carry the week-specific values that shift down.
$shiftCols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# 1) Cache every column of every existing data row before writing anything,
#    since later rows are sourced from earlier rows (and vice versa for the
#    two appended rows).
$allCols = 1..18
$old = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @{}
    foreach ($c in $allCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $old[$r] = $rowVals
}

# 2) Append two new rows (217, 218) that are exact copies of the two rows
#    that fall off the end of the block (old rows 215, 216).
for ($offset = 0; $offset -le 1; $offset++) {
    $srcRow = $lastDataRow - 1 + $offset   # 215, then 216
    $dstRow = $lastDataRow + 1 + $offset   # 217, then 218
    foreach ($c in $allCols) {
        $ws.Cells.Item($dstRow, $c).Value2 = $old[$srcRow][$c]
    }
    $ws.Cells.Item($dstRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# 3) Shift the week-specific columns down by 2 rows for every row from the
#    bottom of the block back up to row 95 (row r takes what used to be in
#    row r-2).
for ($r = $lastDataRow; $r -ge ($firstDataRow + 2); $r--) {
    $srcRow = $r - 2
    foreach ($c in $shiftCols) {
        $ws.Cells.Item($r, $c).Value2 = $old[$srcRow][$c]
    }
}

# 4) The first two rows (93, 94) become the brand-new week's data.
$ws.Cells.Item(93, 4).Value2 = 44557
$ws.Cells.Item(93, 10).Value2 = 1200
$ws.Cells.Item(93, 11).Value2 = 300
$ws.Cells.Item(93, 12).Value2 = 350
$ws.Cells.Item(93, 13).Value2 = 325
$ws.Cells.Item(93, 16).Value2 = 81

$ws.Cells.Item(94, 4).Value2 = 44557
$ws.Cells.Item(94, 10).Value2 = 1200
$ws.Cells.Item(94, 11).Value2 = 300
$ws.Cells.Item(94, 12).Value2 = 350
$ws.Cells.Item(94, 13).Value2 = 325
$ws.Cells.Item(94, 16).Value2 = 65
